$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3957
$ws.Range("I3").Value = 4099
$ws.Range("C4").Value = 1808
$ws.Range("I4").Value = 954
$ws.Range("I5").Value = 377
$ws.Range("I6").Value = 4565
$ws.Range("C7").Value = 28351
$ws.Range("I7").Value = 13952

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I3").Value = 8
$ws.Range("I4").Value = 53
$ws.Range("I7").Value = 441
$ws.Range("I8").Value = 841
$ws.Range("I9").Value = 62
$ws.Range("I11").Value = 213
$ws.Range("I18").Value = 96
$ws.Range("I19").Value = 384
$ws.Range("I20").Value = 344
$ws.Range("I29").Value = 904
$ws.Range("I33").Value = 635
$ws.Range("I36").Value = 195
$ws.Range("I37").Value = 450
$ws.Range("I42").Value = 482
$ws.Range("I44").Value = 99
$ws.Range("I47").Value = 97
$ws.Range("I48").Value = 185
$ws.Range("I49").Value = 117
$ws.Range("I50").Value = 58
$ws.Range("I52").Value = 303
$ws.Range("I54").Value = 319
$ws.Range("I55").Value = 154
$ws.Range("I59").Value = 30
$ws.Range("C63").Value = 241
$ws.Range("I63").Value = 47
$ws.Range("I64").Value = 123
$ws.Range("I65").Value = 314
$ws.Range("I67").Value = 541
$ws.Range("I76").Value = 211
$ws.Range("I79").Value = 378
$ws.Range("I83").Value = 283
$ws.Range("I85").Value = 631
$ws.Range("I87").Value = 26
$ws.Range("I89").Value = 156
$ws.Range("I90").Value = 171
$ws.Range("I92").Value = 42
$ws.Range("I93").Value = 81
$ws.Range("I95").Value = 223
$ws.Range("I96").Value = 151
$ws.Range("I97").Value = 103
$ws.Range("I98").Value = 93
$ws.Range("I99").Value = 260
$ws.Range("C101").Value = 28351
$ws.Range("I101").Value = 13952

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 164
$ws.Range("I3").Value = 255
$ws.Range("I7").Value = 631

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 86
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 96
$ws.Range("I7").Value = 213

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 262
$ws.Range("I4").Value = 47
$ws.Range("I5").Value = 26
$ws.Range("I6").Value = 272
$ws.Range("I7").Value = 841

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 152
$ws.Range("I3").Value = 136
$ws.Range("I7").Value = 441

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 142
$ws.Range("I3").Value = 143
$ws.Range("I4").Value = 30
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 450

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 194
$ws.Range("I4").Value = 27
$ws.Range("I7").Value = 541

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 102
$ws.Range("I3").Value = 88
$ws.Range("I7").Value = 314

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 99
$ws.Range("I3").Value = 109
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 283

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 235
$ws.Range("I7").Value = 635

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I2").Value = 21
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 74
$ws.Range("I3").Value = 62
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 319

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 266
$ws.Range("I3").Value = 310
$ws.Range("I6").Value = 247
$ws.Range("I7").Value = 904

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 384

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 131
$ws.Range("I7").Value = 482

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 48
$ws.Range("I3").Value = 44
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 110
$ws.Range("I3").Value = 120
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 378

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 96
$ws.Range("I6").Value = 110
$ws.Range("I7").Value = 344

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 81

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I2").Value = 24
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 58
$ws.Range("I7").Value = 171

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 8

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I3").Value = 7
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 26
